# Add daily power records
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$tbl = $ws.ListObjects.Item(1)

# Fill in the missing "End Time" value for the existing last row (row 36)
$ws.Range("C36").Value = 0

# Expand the table by one row (this grows the table ref, autoFilter ref,
# and the worksheet dimension to A1:F37)
$newRow = $tbl.ListRows.Add()

# Populate the new row (row 37) with the new daily power record
$ws.Range("A37").Value = 43361
$ws.Range("B37").Value = 0
$ws.Range("D37").Formula = "=(C37-B37)* 1440"
$ws.Range("E37").Formula = "=IF(C37>B37, (C37-B37)*1440, (B37-C37)*1440)"
$ws.Range("F37").Formula = "=ABS((C37-B37)*1440)"

# Update the active selection / scroll position to match the new last row
[void]$ws.Range("C37").Select()
